$d = $word.ActiveDocument

$find = "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός του Περσέα: 16-25 Ιανουαρίου, 7-16 Νοεμβρίου, 6-15 Δεκεμβρίου"
$replace = "2022 Ημερομηνίες παρατήρησης για τον  Αστερισμός του Περσέα: 16-25 Ιανουαρίου, 7-16 Νοεμβρίου, 6-15 Δεκεμβρίου"

$r = $d.Content
$r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
